$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), matching the formatting of the other
# header cells (e.g. G1's bold/centered/bordered style).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

# New data value for the Save column on row 2.
$ws.Range("H2").Value = 1
